$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Cells.Item(6, 8).Value = 12800  # H6: 6288.8887 -> 12800
$ws.Cells.Item(6, 10).Value = 0  # J6: 1080 -> 0
$ws.Cells.Item(6, 12).Value = 0  # L6: 3240 -> 0
$ws.Cells.Item(6, 14).ClearContents()  # N6: -3464 -> (removed)
# Row 12
$ws.Cells.Item(12, 8).Value = 250.5  # H12: 300 -> 250.5
$ws.Cells.Item(12, 9).Value = 250.5  # I12: 300 -> 250.5
$ws.Cells.Item(12, 11).Value = 250.5  # K12: 300 -> 250.5
$ws.Cells.Item(12, 13).Value = -80.5  # M12: -130 -> -80.5
# Row 76
$ws.Cells.Item(76, 8).Value = 6520  # H76: 8333.333000000001 -> 6520
$ws.Cells.Item(76, 10).Value = 5650  # J76: 7500 -> 5650
$ws.Cells.Item(76, 12).Value = 5650  # L76: 7500 -> 5650
$ws.Cells.Item(76, 14).Value = -6280  # N76: -8130 -> -6280
# Row 79
$ws.Cells.Item(79, 8).Value = 6520  # H79: 8333.333000000001 -> 6520
$ws.Cells.Item(79, 10).Value = 5650  # J79: 7500 -> 5650
$ws.Cells.Item(79, 12).Value = 5650  # L79: 7500 -> 5650
$ws.Cells.Item(79, 14).Value = -7834  # N79: -9684 -> -7834
# Row 113
$ws.Cells.Item(113, 8).Value = 11767102  # H113: 11113402 -> 11767102
$ws.Cells.Item(113, 9).Value = 15386909  # I113: 13335554 -> 15386909
$ws.Cells.Item(113, 10).Value = 2727.75  # J113: 2638.6667 -> 2727.75
$ws.Cells.Item(113, 11).Value = 15386909  # K113: 13335554 -> 15386909
$ws.Cells.Item(113, 12).Value = 2727.75  # L113: 2638.6667 -> 2727.75
$ws.Cells.Item(113, 13).Value = -15383655  # M113: -13332300 -> -15383655
$ws.Cells.Item(113, 14).Value = -9235.75  # N113: -9146.6667 -> -9235.75
# Row 129
$ws.Cells.Item(129, 8).Value = 872.9474  # H129: 871.8108 -> 872.9474
$ws.Cells.Item(129, 9).Value = 735.8  # I129: 743.8333 -> 735.8
$ws.Cells.Item(129, 10).Value = 893.7273  # J129: 896.5806 -> 893.7273
$ws.Cells.Item(129, 11).Value = 2207.4  # K129: 2231.4999 -> 2207.4
$ws.Cells.Item(129, 12).Value = 2681.1819  # L129: 2689.7418 -> 2681.1819
$ws.Cells.Item(129, 13).Value = 2792.6  # M129: 2768.5001 -> 2792.6
$ws.Cells.Item(129, 14).Value = -12681.1819  # N129: -12689.7418 -> -12681.1819
# Row 134
$ws.Cells.Item(134, 8).Value = 37797.273  # H134: 38209 -> 37797.273
$ws.Cells.Item(134, 10).Value = 37797.273  # J134: 38209 -> 37797.273
$ws.Cells.Item(134, 12).Value = 37797.273  # L134: 38209 -> 37797.273
$ws.Cells.Item(134, 14).Value = -47937.273  # N134: -48349 -> -47937.273
# Row 137
$ws.Cells.Item(137, 8).Value = 2090.875  # H137: 1858.2778 -> 2090.875
$ws.Cells.Item(137, 9).Value = 1279.6666  # I137: 1131.4783 -> 1279.6666
$ws.Cells.Item(137, 10).Value = 3133.8572  # J137: 3144.1538 -> 3133.8572
$ws.Cells.Item(137, 11).Value = 3838.9998  # K137: 3394.4349 -> 3838.9998
$ws.Cells.Item(137, 12).Value = 9401.571599999999  # L137: 9432.4614 -> 9401.571599999999
$ws.Cells.Item(137, 13).Value = -1288.9998  # M137: -844.4349000000002 -> -1288.9998
$ws.Cells.Item(137, 14).Value = -14501.5716  # N137: -14532.4614 -> -14501.5716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 52632830  # H61: 40001176 -> 52632830
$ws.Cells.Item(61, 9).Value = 66667556  # I61: 50000868 -> 66667556
$ws.Cells.Item(61, 10).Value = 2625  # J61: 2400 -> 2625
$ws.Cells.Item(61, 11).Value = 66667556  # K61: 50000868 -> 66667556
$ws.Cells.Item(61, 12).Value = 2625  # L61: 2400 -> 2625
$ws.Cells.Item(61, 13).Value = -66667344  # M61: -50000656 -> -66667344
$ws.Cells.Item(61, 14).Value = -3049  # N61: -2824 -> -3049
# Row 63
$ws.Cells.Item(63, 8).Value = 2341.5  # H63: 2380.3125 -> 2341.5
$ws.Cells.Item(63, 9).Value = 2183.5454  # I63: 2254 -> 2183.5454
$ws.Cells.Item(63, 10).Value = 2689  # J63: 2590.8333 -> 2689
$ws.Cells.Item(63, 11).Value = 2183.5454  # K63: 2254 -> 2183.5454
$ws.Cells.Item(63, 12).Value = 2689  # L63: 2590.8333 -> 2689
$ws.Cells.Item(63, 13).Value = -1497.5454  # M63: -1568 -> -1497.5454
$ws.Cells.Item(63, 14).Value = -4061  # N63: -3962.8333 -> -4061
# Row 66
$ws.Cells.Item(66, 8).Value = 2341.5  # H66: 2380.3125 -> 2341.5
$ws.Cells.Item(66, 9).Value = 2183.5454  # I66: 2254 -> 2183.5454
$ws.Cells.Item(66, 10).Value = 2689  # J66: 2590.8333 -> 2689
$ws.Cells.Item(66, 11).Value = 10917.727  # K66: 11270 -> 10917.727
$ws.Cells.Item(66, 12).Value = 13445  # L66: 12954.1665 -> 13445
$ws.Cells.Item(66, 13).Value = -7485.726999999999  # M66: -7838 -> -7485.726999999999
$ws.Cells.Item(66, 14).Value = -20309  # N66: -19818.1665 -> -20309
# Row 118
$ws.Cells.Item(118, 8).Value = 33500  # H118: 0 -> 33500
$ws.Cells.Item(118, 9).Value = 20000  # I118: 0 -> 20000
$ws.Cells.Item(118, 10).Value = 38000  # J118: 0 -> 38000
$ws.Cells.Item(118, 11).Value = 20000  # K118: 0 -> 20000
$ws.Cells.Item(118, 12).Value = 38000  # L118: 0 -> 38000
$ws.Cells.Item(118, 13).Value = -18343  # M118: None -> -18343
$ws.Cells.Item(118, 14).Value = -41314  # N118: None -> -41314
# Row 132
$ws.Cells.Item(132, 8).Value = 2608.8076  # H132: 2571.849 -> 2608.8076
$ws.Cells.Item(132, 9).Value = 2021.2  # I132: 1983.1111 -> 2021.2
$ws.Cells.Item(132, 11).Value = 6063.6  # K132: 5949.3333 -> 6063.6
$ws.Cells.Item(132, 13).Value = -3533.6  # M132: -3419.3333 -> -3533.6
# Row 136
$ws.Cells.Item(136, 8).Value = 52632830  # H136: 40001176 -> 52632830
$ws.Cells.Item(136, 9).Value = 66667556  # I136: 50000868 -> 66667556
$ws.Cells.Item(136, 10).Value = 2625  # J136: 2400 -> 2625
$ws.Cells.Item(136, 11).Value = 200002668  # K136: 150002604 -> 200002668
$ws.Cells.Item(136, 12).Value = 7875  # L136: 7200 -> 7875
$ws.Cells.Item(136, 13).Value = -200000118  # M136: -150000054 -> -200000118
$ws.Cells.Item(136, 14).Value = -12975  # N136: -12300 -> -12975

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Cells.Item(105, 8).Value = 201978720  # H105: 144270700 -> 201978720
$ws.Cells.Item(105, 9).Value = 336630300  # I105: 168315490 -> 336630300
$ws.Cells.Item(105, 10).Value = 1375  # J105: 2000 -> 1375
$ws.Cells.Item(105, 11).Value = 336630300  # K105: 168315490 -> 336630300
$ws.Cells.Item(105, 12).Value = 1375  # L105: 2000 -> 1375
$ws.Cells.Item(105, 13).Value = -336628553  # M105: -168313743 -> -336628553
$ws.Cells.Item(105, 14).Value = -4869  # N105: -5494 -> -4869
# Row 116
$ws.Cells.Item(116, 8).Value = 0  # H116: 30742 -> 0
$ws.Cells.Item(116, 10).Value = 0  # J116: 30742 -> 0
$ws.Cells.Item(116, 12).ClearContents()  # L116: 30742 -> (removed)
$ws.Cells.Item(116, 14).Value = 0  # N116: -39920 -> 0
# Row 134
$ws.Cells.Item(134, 8).Value = 3514.4102  # H134: 3441.55 -> 3514.4102
$ws.Cells.Item(134, 10).Value = 9646.083000000001  # J134: 8950.23 -> 9646.083000000001
$ws.Cells.Item(134, 12).Value = 28938.249  # L134: 26850.69 -> 28938.249
$ws.Cells.Item(134, 14).Value = -34008.249  # N134: -31920.69 -> -34008.249

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Cells.Item(62, 8).Value = 13335508  # H62: 12502220 -> 13335508
$ws.Cells.Item(62, 9).Value = 2293.8462  # I62: 2352 -> 2293.8462
$ws.Cells.Item(62, 10).Value = 100001400  # J62: 33335334 -> 100001400
$ws.Cells.Item(62, 11).Value = 2293.8462  # K62: 2352 -> 2293.8462
$ws.Cells.Item(62, 12).Value = 100001400  # L62: 33335334 -> 100001400
$ws.Cells.Item(62, 13).Value = -1669.8462  # M62: -1728 -> -1669.8462
$ws.Cells.Item(62, 14).Value = -100002648  # N62: -33336582 -> -100002648
# Row 65
$ws.Cells.Item(65, 8).Value = 13335508  # H65: 12502220 -> 13335508
$ws.Cells.Item(65, 9).Value = 2293.8462  # I65: 2352 -> 2293.8462
$ws.Cells.Item(65, 10).Value = 100001400  # J65: 33335334 -> 100001400
$ws.Cells.Item(65, 11).Value = 11469.231  # K65: 11760 -> 11469.231
$ws.Cells.Item(65, 12).Value = 500007000  # L65: 166676670 -> 500007000
$ws.Cells.Item(65, 13).Value = -8349.231  # M65: -8640 -> -8349.231
$ws.Cells.Item(65, 14).Value = -500013240  # N65: -166682910 -> -500013240
# Row 132
$ws.Cells.Item(132, 8).Value = 1524.975  # H132: 1411.341 -> 1524.975
$ws.Cells.Item(132, 9).Value = 1196.1923  # I132: 1073.3667 -> 1196.1923
$ws.Cells.Item(132, 11).Value = 3588.5769  # K132: 3220.1001 -> 3588.5769
$ws.Cells.Item(132, 13).Value = -1058.5769  # M132: -690.1001000000001 -> -1058.5769

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 45003520  # H70: 40912656 -> 45003520
$ws.Cells.Item(70, 9).Value = 35717930  # I70: 31253688 -> 35717930
$ws.Cells.Item(70, 11).Value = 35717930  # K70: 31253688 -> 35717930
$ws.Cells.Item(70, 13).Value = -35717660  # M70: -31253418 -> -35717660
# Row 73
$ws.Cells.Item(73, 8).Value = 45003520  # H73: 40912656 -> 45003520
$ws.Cells.Item(73, 9).Value = 35717930  # I73: 31253688 -> 35717930
$ws.Cells.Item(73, 11).Value = 35717930  # K73: 31253688 -> 35717930
$ws.Cells.Item(73, 13).Value = -35716994  # M73: -31252752 -> -35716994
# Row 80
$ws.Cells.Item(80, 8).Value = 3963.5454  # H80: 3288.8823 -> 3963.5454
$ws.Cells.Item(80, 9).Value = 1800  # I80: 1952 -> 1800
$ws.Cells.Item(80, 10).Value = 4179.9  # J80: 4018.0908 -> 4179.9
$ws.Cells.Item(80, 11).Value = 1800  # K80: 1952 -> 1800
$ws.Cells.Item(80, 12).Value = 4179.9  # L80: 4018.0908 -> 4179.9
$ws.Cells.Item(80, 13).Value = -802  # M80: -954 -> -802
$ws.Cells.Item(80, 14).Value = -6175.9  # N80: -6014.0908 -> -6175.9
# Row 83
$ws.Cells.Item(83, 8).Value = 3963.5454  # H83: 3288.8823 -> 3963.5454
$ws.Cells.Item(83, 9).Value = 1800  # I83: 1952 -> 1800
$ws.Cells.Item(83, 10).Value = 4179.9  # J83: 4018.0908 -> 4179.9
$ws.Cells.Item(83, 11).Value = 9000  # K83: 9760 -> 9000
$ws.Cells.Item(83, 12).Value = 20899.5  # L83: 20090.454 -> 20899.5
$ws.Cells.Item(83, 13).Value = -4008  # M83: -4768 -> -4008
$ws.Cells.Item(83, 14).Value = -30883.5  # N83: -30074.454 -> -30883.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 114
$ws.Cells.Item(114, 8).Value = 48000  # H114: 47249.25 -> 48000
$ws.Cells.Item(114, 10).Value = 48000  # J114: 47249.25 -> 48000
$ws.Cells.Item(114, 12).Value = 48000  # L114: 47249.25 -> 48000
$ws.Cells.Item(114, 14).Value = -56678  # N114: -55927.25 -> -56678
# Row 136
$ws.Cells.Item(136, 8).Value = 1644.409  # H136: 1652.591 -> 1644.409
$ws.Cells.Item(136, 9).Value = 1515.8889  # I136: 1525.8889 -> 1515.8889
$ws.Cells.Item(136, 11).Value = 4547.6667  # K136: 4577.6667 -> 4547.6667
$ws.Cells.Item(136, 13).Value = -1997.6667  # M136: -2027.6667 -> -1997.6667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Cells.Item(46, 8).Value = 43999.25  # H46: 45998.75 -> 43999.25
$ws.Cells.Item(46, 10).Value = 43999.25  # J46: 45998.75 -> 43999.25
$ws.Cells.Item(46, 12).Value = 43999.25  # L46: 45998.75 -> 43999.25
$ws.Cells.Item(46, 14).Value = -44461.25  # N46: -46460.75 -> -44461.25
# Row 93
$ws.Cells.Item(93, 8).Value = 8425  # H93: 8416 -> 8425
$ws.Cells.Item(93, 10).Value = 8425  # J93: 8416 -> 8425
$ws.Cells.Item(93, 12).Value = 8425  # L93: 8416 -> 8425
$ws.Cells.Item(93, 14).Value = -13417  # N93: -13408 -> -13417
# Row 132
$ws.Cells.Item(132, 8).Value = 3080.6736  # H132: 3361.4783 -> 3080.6736
$ws.Cells.Item(132, 9).Value = 3111.5854  # I132: 3465.2632 -> 3111.5854
$ws.Cells.Item(132, 10).Value = 2922.25  # J132: 2868.5 -> 2922.25
$ws.Cells.Item(132, 11).Value = 9334.7562  # K132: 10395.7896 -> 9334.7562
$ws.Cells.Item(132, 12).Value = 8766.75  # L132: 8605.5 -> 8766.75
$ws.Cells.Item(132, 13).Value = -6804.7562  # M132: -7865.7896 -> -6804.7562
$ws.Cells.Item(132, 14).Value = -13826.75  # N132: -13665.5 -> -13826.75
# Row 134
$ws.Cells.Item(134, 8).Value = 43999.25  # H134: 45998.75 -> 43999.25
$ws.Cells.Item(134, 10).Value = 43999.25  # J134: 45998.75 -> 43999.25
$ws.Cells.Item(134, 12).Value = 131997.75  # L134: 137996.25 -> 131997.75
$ws.Cells.Item(134, 14).Value = -137067.75  # N134: -143066.25 -> -137067.75
# Row 135
$ws.Cells.Item(135, 8).Value = 72495  # H135: 64996.668 -> 72495
$ws.Cells.Item(135, 10).Value = 72495  # J135: 64996.668 -> 72495
$ws.Cells.Item(135, 12).Value = 72495  # L135: 64996.668 -> 72495
$ws.Cells.Item(135, 14).Value = -82635  # N135: -75136.66800000001 -> -82635
# Row 136
$ws.Cells.Item(136, 8).Value = 1432.7273  # H136: 1587.2222 -> 1432.7273
$ws.Cells.Item(136, 9).Value = 1084.4445  # I136: 1183.5714 -> 1084.4445
$ws.Cells.Item(136, 11).Value = 3253.3335  # K136: 3550.7142 -> 3253.3335
$ws.Cells.Item(136, 13).Value = -703.3335000000002  # M136: -1000.7142 -> -703.3335000000002
# Row 141
$ws.Cells.Item(141, 8).Value = 50612.777  # H141: 48923 -> 50612.777
$ws.Cells.Item(141, 10).Value = 50612.777  # J141: 48923 -> 50612.777
$ws.Cells.Item(141, 12).Value = 50612.777  # L141: 48923 -> 50612.777
$ws.Cells.Item(141, 14).Value = -60972.777  # N141: -59283 -> -60972.777
